$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended to the table (rows 60-69), matching the
# existing table's column layout (A: subject id [text], B: trial,
# C: start, D: end).
$newRows = @(
    @("40", 3, 197, 198),
    @("112", 4, 15, 21),
    @("112", 5, 15, 21),
    @("112", 6, 15, 21),
    @("112", 7, 15, 21),
    @("112", 8, 15, 21),
    @("112", 9, 15, 21),
    @("112", 10, 15, 21),
    @("112", 11, 15, 21),
    @("112", 12, 15, 21)
)

$startRow = 60
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Scroll the viewport and update the selection to reflect where the
# editor left off (best-effort; not all hosts expose window scrolling).
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 49
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("E68").Select()

$wb.Save()
